$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.MoveEnd(1, -1)
    $r.Text = $newText
}

function Replace-Text($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $newText, 2) | Out-Null
}

# --- Row ID=1 (table row 2): Priority 1->4, Sprint 1->4 ---
Set-CellText $t 2 5 "4"
Set-CellText $t 2 6 "4"

# --- Row ID=2 (table row 3): Priority 2->5, Sprint 2->4 ---
Set-CellText $t 3 5 "5"
Set-CellText $t 3 6 "4"

# --- Row ID=3 (table row 4) ---
# "I want" cell: "see specific product details" -> "Add products details"
Replace-Text "see specific product details" "Add products details"

# "So that" cell: " can give detail information of the product to the customers"
#              -> " can let the customer know detail information of the product."
Replace-Text " can give detail information of the product to the customers" `
             " can let the customer know detail information of the product."

# Priority 2->3, Sprint 2->3, Status "To be started" -> "Done"
Set-CellText $t 4 5 "3"
Set-CellText $t 4 6 "3"
Set-CellText $t 4 7 "Done"

# --- Row ID=4 (table row 5): Priority 3->1, Sprint 3->1, Status -> Done ---
Set-CellText $t 5 5 "1"
Set-CellText $t 5 6 "1"
Set-CellText $t 5 7 "Done"

# --- Row ID=5 (table row 6) ---
# "I want" cell: "search specific customer's record" -> "search customer's record"
Replace-Text "search specific customer" "search customer"

# "So that" cell: "The customer's regularity of products can be checked and kept updated"
#              -> "The customer's details can be found in time of necessity"
$apos = [char]0x2019
Replace-Text "The customer$($apos)s regularity of products can be checked and kept updated" `
             "The customer$($apos)s details can be found in time of necessity"

# Priority 3->2, Sprint 3->2, Status -> Done
Set-CellText $t 6 5 "2"
Set-CellText $t 6 6 "2"
Set-CellText $t 6 7 "Done"

# --- Row ID=6 (table row 7): Priority 4->6 ---
Set-CellText $t 7 5 "6"

# --- Row ID=7 (table row 8): Priority 4->7 ---
Set-CellText $t 8 5 "7"

# --- Row ID=8 (table row 9): Priority 4->8 ---
Set-CellText $t 9 5 "8"

# --- Row ID=9 (table row 10): Priority 5->9 ---
Set-CellText $t 10 5 "9"
